$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename shared string "Neutrophils" -> "Resolving-Mac" for target cluster cells (D6, D11)
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("D11").Value = "Resolving-Mac"

# Row 2
$ws.Range("I2").Value = 0.02328126719340038
$ws.Range("J2").Value = 0.02328126719340038
$ws.Range("M2").Value = 2.157506
$ws.Range("N2").Value = 6.472517999999999
$ws.Range("O2").Value = 0.3549648016839517
$ws.Range("P2").Value = 0.3549648016839516
$ws.Range("Q2").Value = 0.04549532902199999
$ws.Range("R2").Value = 0.4094579611979999
$ws.Range("S2").Value = 0.008264030392256457
$ws.Range("T2").Value = 0.008264030392256457

# Row 3
$ws.Range("I3").Value = 0.02328126719340038
$ws.Range("J3").Value = 0.02328126719340038
$ws.Range("O3").Value = 0.4793705560628122
$ws.Range("P3").Value = 0.4793705560628121
$ws.Range("S3").Value = 0.01116035400034725
$ws.Range("T3").Value = 0.01116035400034725

# Row 4
$ws.Range("I4").Value = 0.02328126719340038
$ws.Range("J4").Value = 0.02328126719340038
$ws.Range("M4").Value = 0.018986
$ws.Range("N4").Value = 0.05695799999999999
$ws.Range("O4").Value = 0.003123681567871193
$ws.Range("P4").Value = 0.003123681567871192
$ws.Range("Q4").Value = 0.000400357782
$ws.Range("R4").Value = 0.003603220038
$ws.Range("S4").Value = 0.00007272326520870908
$ws.Range("T4").Value = 0.00007272326520870907

# Row 5
$ws.Range("I5").Value = 0.02328126719340038
$ws.Range("J5").Value = 0.02328126719340038
$ws.Range("M5").Value = 0.9848966666666668
$ws.Range("N5").Value = 2.95469
$ws.Range("O5").Value = 0.1620406385718132
$ws.Range("P5").Value = 0.1620406385718132
$ws.Range("Q5").Value = 0.02076851601
$ws.Range("R5").Value = 0.18691664409
$ws.Range("S5").Value = 0.003772511402779604
$ws.Range("T5").Value = 0.003772511402779603

# Row 6
$ws.Range("I6").Value = 0.02328126719340038
$ws.Range("J6").Value = 0.02328126719340038
$ws.Range("M6").Value = 0.003041
$ws.Range("N6").Value = 0.009122999999999999
$ws.Range("O6").Value = 0.0005003221135518961
$ws.Range("P6").Value = 0.000500322113551896
$ws.Range("Q6").Value = 0.000064125567
$ws.Range("R6").Value = 0.0005771301029999999
$ws.Range("S6").Value = 0.0000116481328083685
$ws.Range("T6").Value = 0.0000116481328083685

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.8846626666666667
$ws.Range("H7").Value = 2.653988
$ws.Range("I7").Value = 0.9767187328065996
$ws.Range("J7").Value = 0.9767187328065997
$ws.Range("M7").Value = 2.157506
$ws.Range("N7").Value = 6.472517999999999
$ws.Range("O7").Value = 0.3549648016839517
$ws.Range("P7").Value = 0.3549648016839516
$ws.Range("Q7").Value = 1.908665011309333
$ws.Range("R7").Value = 17.177985101784
$ws.Range("S7").Value = 0.3467007712916952
$ws.Range("T7").Value = 0.3467007712916952

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.8846626666666667
$ws.Range("H8").Value = 2.653988
$ws.Range("I8").Value = 0.9767187328065996
$ws.Range("J8").Value = 0.9767187328065997
$ws.Range("O8").Value = 0.4793705560628122
$ws.Range("P8").Value = 0.4793705560628121
$ws.Range("Q8").Value = 2.577601507159111
$ws.Range("R8").Value = 23.198413564432
$ws.Range("S8").Value = 0.4682102020624649
$ws.Range("T8").Value = 0.4682102020624649

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.8846626666666667
$ws.Range("H9").Value = 2.653988
$ws.Range("I9").Value = 0.9767187328065996
$ws.Range("J9").Value = 0.9767187328065997
$ws.Range("M9").Value = 0.018986
$ws.Range("N9").Value = 0.05695799999999999
$ws.Range("O9").Value = 0.003123681567871193
$ws.Range("P9").Value = 0.003123681567871192
$ws.Range("Q9").Value = 0.01679620538933333
$ws.Range("R9").Value = 0.151165848504
$ws.Range("S9").Value = 0.003050958302662484
$ws.Range("T9").Value = 0.003050958302662483

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.8846626666666667
$ws.Range("H10").Value = 2.653988
$ws.Range("I10").Value = 0.9767187328065996
$ws.Range("J10").Value = 0.9767187328065997
$ws.Range("M10").Value = 0.9848966666666668
$ws.Range("N10").Value = 2.95469
$ws.Range("O10").Value = 0.1620406385718132
$ws.Range("P10").Value = 0.1620406385718132
$ws.Range("Q10").Value = 0.8713013115244446
$ws.Range("R10").Value = 7.841711803720001
$ws.Range("S10").Value = 0.1582681271690336
$ws.Range("T10").Value = 0.1582681271690336

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.8846626666666667
$ws.Range("H11").Value = 2.653988
$ws.Range("I11").Value = 0.9767187328065996
$ws.Range("J11").Value = 0.9767187328065997
$ws.Range("M11").Value = 0.003041
$ws.Range("N11").Value = 0.009122999999999999
$ws.Range("O11").Value = 0.0005003221135518961
$ws.Range("P11").Value = 0.000500322113551896
$ws.Range("Q11").Value = 0.002690259169333333
$ws.Range("R11").Value = 0.024212332524
$ws.Range("S11").Value = 0.0004886739807435275
$ws.Range("T11").Value = 0.0004886739807435275
